$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a brand-new "Pytanie czwarte" question block right before
#    the "Pytanie piąte" paragraph.
# ------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Trim() -eq "Pytanie piąte") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    $newLines = @(
        "Pytanie czwarte",
        "○ Odpowiedź czwarta",
        "○ Odpowiedź druga",
        "○ Odpowiedź trzecia",
        "○ Odpowiedź pierwsza",
        "○ Odpowiedź piąta"
    )
    $target = $d.Paragraphs($targetIndex)
    for ($i = $newLines.Length - 1; $i -ge 0; $i--) {
        $target.Range.InsertBefore("$($newLines[$i])`r")
    }
}

# ------------------------------------------------------------------
# 2) Rename "Pytanie drugie" -> "Pytanie trzecie".
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Trim() -eq "Pytanie drugie") {
        $p.Range.Text = "Pytanie trzecie"
        break
    }
}

# ------------------------------------------------------------------
# 3) The renamed "Pytanie trzecie" block now has an extra trailing
#    "○ Odpowiedź czwarta" answer that needs to be removed (the block
#    should only keep pierwsza/druga/trzecia).
# ------------------------------------------------------------------
$questionIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Trim() -eq "Pytanie trzecie") {
        $questionIndex = $i
        break
    }
}

if ($questionIndex -ge 1) {
    for ($i = $questionIndex + 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        $text = $p.Range.Text.Trim()
        if ($text -eq "○ Odpowiedź czwarta") {
            $p.Range.Delete()
            break
        }
        if ($text -eq "" -or $text.StartsWith("Pytanie")) {
            break
        }
    }
}

# ------------------------------------------------------------------
# 4) Header: "Test 1" -> "Test 2".
# ------------------------------------------------------------------
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections($s)
    for ($h = 1; $h -le 3; $h++) {
        $hdr = $section.Headers($h)
        if ($hdr.Exists) {
            for ($i = 1; $i -le $hdr.Range.Paragraphs.Count; $i++) {
                $p = $hdr.Range.Paragraphs($i)
                if ($p.Range.Text.Trim() -eq "Test 1") {
                    $p.Range.Text = "Test 2"
                }
            }
        }
    }
}
